$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current blank separator row (row 74),
# shifting the blank separator row and the summary rows down by two.
$ws.Rows("74:75").Insert()

# Row 74: new data row (2014-03-17, 15:30 - 17:30)
$ws.Range("A74").Value = 2014
$ws.Range("B74").Value = 3
$ws.Range("C74").Value = 17
$ws.Range("D74").Value = 0.64583333333333337
$ws.Range("E74").Value = 0.72916666666666663
$ws.Range("F74").Formula = "=(E74-D74)*24*60"
$ws.Range("G74").Formula = "=F74/60"

# Row 75: new data row (2014-03-17, 20:30 - 22:00)
$ws.Range("A75").Value = 2014
$ws.Range("B75").Value = 3
$ws.Range("C75").Value = 17
$ws.Range("D75").Value = 0.85416666666666663
$ws.Range("E75").Value = 0.91666666666666663
$ws.Range("F75").Formula = "=(E75-D75)*24*60"
$ws.Range("G75").Formula = "=F75/60"

# Update the selected cell to match the new layout.
$null = $ws.Range("A76").Select()
